# Add a new "Walmart Inc" row to the ticker list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the new row. The RIC / Yahoo Ticker are entered before the
# Company name (matching the order the new strings were appended to the
# workbook's shared string table), and the Unit last.
$ws.Cells.Item(8, 2).Value = "WMT.N"
$ws.Cells.Item(8, 3).Value = "WMT"
$ws.Cells.Item(8, 1).Value = "Walmart Inc"
$ws.Cells.Item(8, 4).Value = "USD"

# Copy the "Unit" cell formatting from the row above (D7) so the new
# D8 cell matches the style used by the rest of the Unit column.
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Leave the active cell on the newly added row, column A.
[void]$ws.Range("A8").Select()
